$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.452.93'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.617.51'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue $ws.Range('D5') '212.59'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D8') '0.0608'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D9') '0.245'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +1.45%  '
Set-TextValue $ws.Range('D11') '0.0847'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '1.844.61'
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').Value = '1.614.24'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('E15').Value = '  +0.05%  '
Set-TextValue $ws.Range('D16') '63.83'
$ws.Range('E16').Value = '  +0.07%  '
Set-TextValue $ws.Range('D17') '237.59'
$ws.Range('E17').Value = '  +10.15%  '
$ws.Range('D18').Value = '26.448.16'
$ws.Range('E18').Value = '  +0.72%  '
Set-TextValue $ws.Range('D19') '7.77'
$ws.Range('E19').Value = '  +5.76%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +0.01%  '
Set-TextValue $ws.Range('D22') '4.30'
$ws.Range('E22').Value = '  +0.06%  '
Set-TextValue $ws.Range('D23') '2.19'
$ws.Range('E23').Value = '  +4.29%  '
$ws.Range('E24').Value = '  +0.53%  '
Set-TextValue $ws.Range('D25') '147.09'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  -0.06%  '
Set-TextValue $ws.Range('D27') '7.02'
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('E28').Value = '  +0.03%  '
Set-TextValue $ws.Range('D29') '15.51'
$ws.Range('E29').Value = '  +2.69%  '
Set-TextValue $ws.Range('D30') '0.0497'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '1.516.92'
$ws.Range('E32').Value = '  +6.84%  '
Set-TextValue $ws.Range('D33') '3.24'
$ws.Range('E33').Value = '  +1.48%  '
Set-TextValue $ws.Range('D34') '2.97'
$ws.Range('E34').Value = '  +0.25%  '
Set-TextValue $ws.Range('D35') '1.54'
$ws.Range('E35').Value = '  +5.64%  '
$ws.Range('E36').Value = '  +0.05%  '
Set-TextValue $ws.Range('D37') '0.567'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('E38').Value = '  +0.26%  '
Set-TextValue $ws.Range('D39') '0.830'
$ws.Range('E39').Value = '  +0.48%  '
Set-TextValue $ws.Range('D40') '5.92'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('D43').Value = '1.755.91'
$ws.Range('E43').Value = '  +1.55%  '
Set-TextValue $ws.Range('D44') '0.760'
$ws.Range('E44').Value = '  -0.20%  '
Set-TextValue $ws.Range('D45') '0.914'
$ws.Range('E45').Value = '  -3.04%  '
Set-TextValue $ws.Range('D46') '61.48'
$ws.Range('E46').Value = '  +0.96%  '
Set-TextValue $ws.Range('D47') '90.19'
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D48') '1.50'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D49') '0.0502'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D50') '0.0961'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '7.47'
$ws.Range('E51').Value = '  +0.49%  '
